$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 10
$ws.Range("S10").Value = 4.3
$ws.Range("W10").Value = 1.57
$ws.Range("X10").Value = 2.25
$ws.Range("AA10").Value = 6
$ws.Range("AE10").Value = 21
$ws.Range("AG10").Value = 6.5

# Row 25
$ws.Range("G25").Value = 1.95
$ws.Range("I25").Value = 4.5
$ws.Range("M25").Value = 1.1
$ws.Range("N25").Value = 7
$ws.Range("Q25").Value = 2.5
$ws.Range("R25").Value = 1.5
$ws.Range("AB25").Value = 8
$ws.Range("AD25").Value = 17
$ws.Range("AG25").Value = 6.5
$ws.Range("AI25").Value = 19
$ws.Range("AL25").Value = 9
$ws.Range("AR25").Value = 1.93
$ws.Range("AS25").Value = 1.93

# Row 60
$ws.Range("G60").Value = 1.55
$ws.Range("H60").Value = 4.2
$ws.Range("I60").Value = 5.5
$ws.Range("J60").Value = 2.1
$ws.Range("L60").Value = 5.5
$ws.Range("O60").Value = 1.18
$ws.Range("P60").Value = 4.5
$ws.Range("Q60").Value = 1.62
$ws.Range("R60").Value = 2.25
$ws.Range("S60").Value = 2.05
$ws.Range("T60").Value = 1.8
$ws.Range("U60").Value = 2.5
$ws.Range("V60").Value = 1.5
$ws.Range("Y60").Value = 1.7
$ws.Range("AD60").Value = 12

# Row 61
$ws.Range("G61").Value = 1.37
$ws.Range("H61").Value = 4.75
$ws.Range("I61").Value = 6.7
$ws.Range("J61").Value = 1.8
$ws.Range("K61").Value = 2.55
$ws.Range("L61").Value = 5.8
$ws.Range("Y61").Value = 1.62
$ws.Range("Z61").Value = 2.02
$ws.Range("AA61").Value = 9.75
$ws.Range("AB61").Value = 8
$ws.Range("AD61").Value = 10
$ws.Range("AF61").Value = 20
$ws.Range("AG61").Value = 18
$ws.Range("AH61").Value = 10
$ws.Range("AI61").Value = 16.5
$ws.Range("AJ61").Value = 55
$ws.Range("AK61").Value = 350
$ws.Range("AM61").Value = 50
$ws.Range("AN61").Value = 21
$ws.Range("AP61").Value = 60
$ws.Range("AQ61").Value = 50

# Row 71
$ws.Range("I71").Value = 11
$ws.Range("M71").Value = 1.03
$ws.Range("N71").Value = 15
$ws.Range("Q71").Value = 1.62
$ws.Range("R71").Value = 2.25
$ws.Range("U71").Value = 2.5
$ws.Range("V71").Value = 1.5
$ws.Range("AO71").Value = 126
